# hosadatbazis.xlsx - KepAdatbazis kiegeszitese
# Replace the old "hero picture" path strings with the newly uploaded
# "enemy / golem" picture paths, and update each sheet's view/selection
# state to reflect the new active selections.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Harcosok ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C2").Value = "Images\Karakterek\ellenseg.png"
$ws1.Range("C3").Value = "Images\Karakterek\ellenseg2.png"
$ws1.Range("C4").Value = "Images\Karakterek\ellenseg3.png"
$ws1.Range("C5").Value = "Images\Karakterek\golem1.png"
$ws1.Columns.Item(3).ColumnWidth = 39.5

# --- Sheet 2: Ijaszok ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C2").Value = "Images\Karakterek\golem2.png"
$ws2.Range("C3").Value = "Images\Karakterek\golem3.png"
$ws2.Range("C4").Value = "Images\Karakterek\golem4.png"

# --- Sheet 3: Magusok ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C2").Value = "Images\Karakterek\golem2.png"
$ws3.Range("C3").Value = "Images\Karakterek\golem3.png"
$ws3.Range("C4").Value = "Images\Karakterek\golem4.png"

# Update each sheet's own saved selection/view state
$ws2.Activate()
$ws2.Range("C2:C4").Select()

$ws3.Activate()
$ws3.Range("C8").Select()

# Harcosok ends up the active tab, with C6:C8 selected
$ws1.Activate()
$ws1.Range("C6:C8").Select()
